$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells for the Python object-oriented / Eigen columns
$ws.Range("F1").Value = "Object-Oriented Kalman Python"
$ws.Range("G1").Value = "Kalman with Eigen Python (AvgTime for 1000 Runs)"
$ws.Range("F1:G1").Font.Bold = $true

# Column F data (Object-Oriented Kalman Python timings)
$ws.Range("F2").Value = 158.07
$ws.Range("F3").Value = 74.86
$ws.Range("F4").Value = 59.13
$ws.Range("F5").Value = 52.93
$ws.Range("F6").Value = 59.13
$ws.Range("F7").Value = 49.11
$ws.Range("F8").Value = 49.83
$ws.Range("F9").Value = 52.21
$ws.Range("F10").Value = 55.79
$ws.Range("F11").Value = 54.12

# Averages row, including G12 which will error (#DIV/0!) since column G has no data rows
$ws.Range("F12:G12").Formula = "=AVERAGE(F2:F11)"

# Column widths: F matches E (21.5), G matches C/D (21.6640625)
$ws.Columns("F").ColumnWidth = 20.666666666666668
$ws.Columns("G").ColumnWidth = 20.830729166666668

# Active cell moves to G1 after editing
$ws.Range("G1").Select()
